$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header labels: BP1/BQ1 ("average_doctor" <-> "average_doctor_old")
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Update data values (new Harvard case classification stats)
$ws.Range("E4").Value = 0.412
$ws.Range("F4").Value = 0.08
$ws.Range("G4").Value = 0.283
$ws.Range("N4").Value = 0.415
$ws.Range("O4").Value = 0.067
$ws.Range("P4").Value = 0.26
$ws.Range("W4").Value = 0.226
$ws.Range("X4").Value = 0.104
$ws.Range("Y4").Value = 0.322
$ws.Range("AI4").Value = 0.219
$ws.Range("AJ4").Value = 0.067
$ws.Range("AK4").Value = 0.258
$ws.Range("AU4").Value = 0.151
$ws.Range("AV4").Value = 0.028
$ws.Range("AW4").Value = 0.166
$ws.Range("BA4").Value = 1.946
$ws.Range("BB4").Value = 0.161
$ws.Range("BC4").Value = 0.401
$ws.Range("BG4").Value = 0.722
$ws.Range("BH4").Value = 0.14
$ws.Range("BI4").Value = 0.374
$ws.Range("BM4").Value = 0.694
$ws.Range("BN4").Value = 0.083
$ws.Range("BO4").Value = 0.288
$ws.Range("BP4").Value = 0.649
$ws.Range("BQ4").Value = 0.65
$ws.Range("E5").Value = 0.533
$ws.Range("F5").Value = 0.092
$ws.Range("G5").Value = 0.304
$ws.Range("N5").Value = 0.754
$ws.Range("O5").Value = 0.079
$ws.Range("P5").Value = 0.28
$ws.Range("W5").Value = 0.228
$ws.Range("X5").Value = 0.108
$ws.Range("Y5").Value = 0.329
$ws.Range("AI5").Value = 0.261
$ws.Range("AJ5").Value = 0.097
$ws.Range("AK5").Value = 0.311
$ws.Range("AU5").Value = 0.308
$ws.Range("AV5").Value = 0.103
$ws.Range("AW5").Value = 0.32
$ws.Range("BA5").Value = 1.385
$ws.Range("BB5").Value = 0.08799999999999999
$ws.Range("BC5").Value = 0.296
$ws.Range("BG5").Value = 0.411
$ws.Range("BH5").Value = 0.054
$ws.Range("BI5").Value = 0.232
$ws.Range("BM5").Value = 0.582
$ws.Range("BO5").Value = 0.273
$ws.Range("BP5").Value = 0.462
$ws.Range("BQ5").Value = 0.455
$ws.Range("E6").Value = 0.465
$ws.Range("N6").Value = 0.535
$ws.Range("W6").Value = 0.227
$ws.Range("AI6").Value = 0.238
$ws.Range("AU6").Value = 0.203
$ws.Range("BA6").Value = 1.608
$ws.Range("BG6").Value = 0.524
$ws.Range("BM6").Value = 0.633
$ws.Range("BP6").Value = 0.536
$ws.Range("BQ6").Value = 0.532
$ws.Range("E7").Value = 0.503
$ws.Range("N7").Value = 0.648
$ws.Range("W7").Value = 0.228
$ws.Range("AI7").Value = 0.251
$ws.Range("AU7").Value = 0.255
$ws.Range("BA7").Value = 1.465
$ws.Range("BG7").Value = 0.45
$ws.Range("BM7").Value = 0.601
$ws.Range("BP7").Value = 0.488
$ws.Range("BQ7").Value = 0.483
$ws.Range("E8").Value = 0.5679999999999999
$ws.Range("F8").Value = 0.117
$ws.Range("G8").Value = 0.342
$ws.Range("N8").Value = 0.766
$ws.Range("O8").Value = 0.065
$ws.Range("P8").Value = 0.256
$ws.Range("W8").Value = 0.229
$ws.Range("X8").Value = 0.11
$ws.Range("Y8").Value = 0.331
$ws.Range("AI8").Value = 0.239
$ws.Range("AJ8").Value = 0.096
$ws.Range("AK8").Value = 0.311
$ws.Range("AU8").Value = 0.24
$ws.Range("AV8").Value = 0.07199999999999999
$ws.Range("AW8").Value = 0.268
$ws.Range("BA8").Value = 1.702
$ws.Range("BG8").Value = 0.555
$ws.Range("BH8").Value = 0.103
$ws.Range("BI8").Value = 0.321
$ws.Range("BM8").Value = 0.703
$ws.Range("BN8").Value = 0.066
$ws.Range("BO8").Value = 0.257
$ws.Range("BP8").Value = 0.5669999999999999
$ws.Range("BQ8").Value = 0.574
$ws.Range("E9").Value = 0.489
$ws.Range("F9").Value = 0.25
$ws.Range("G9").Value = 0.5
$ws.Range("N9").Value = 0.667
$ws.Range("O9").Value = 0.222
$ws.Range("P9").Value = 0.471
$ws.Range("W9").Value = 0.133
$ws.Range("X9").Value = 0.116
$ws.Range("Y9").Value = 0.34
$ws.Range("AI9").Value = 0.133
$ws.Range("AJ9").Value = 0.116
$ws.Range("AK9").Value = 0.34
$ws.Range("BA9").Value = 1.6
$ws.Range("BB9").Value = 0.24
$ws.Range("BC9").Value = 0.49
$ws.Range("BG9").Value = 0.578
$ws.Range("BH9").Value = 0.244
$ws.Range("BI9").Value = 0.494
$ws.Range("BM9").Value = 0.622
$ws.Range("BN9").Value = 0.235
$ws.Range("BO9").Value = 0.485
$ws.Range("BP9").Value = 0.533
$ws.Range("BQ9").Value = 0.528
$ws.Range("E10").Value = 0.622
$ws.Range("F10").Value = 0.235
$ws.Range("G10").Value = 0.485
$ws.Range("N10").Value = 0.867
$ws.Range("O10").Value = 0.116
$ws.Range("P10").Value = 0.34
$ws.Range("W10").Value = 0.267
$ws.Range("X10").Value = 0.196
$ws.Range("Y10").Value = 0.442
$ws.Range("AI10").Value = 0.267
$ws.Range("AJ10").Value = 0.196
$ws.Range("AK10").Value = 0.442
$ws.Range("AU10").Value = 0.222
$ws.Range("AV10").Value = 0.173
$ws.Range("AW10").Value = 0.416
$ws.Range("BA10").Value = 1.956
$ws.Range("BB10").Value = 0.249
$ws.Range("BC10").Value = 0.499
$ws.Range("BG10").Value = 0.622
$ws.Range("BH10").Value = 0.235
$ws.Range("BI10").Value = 0.485
$ws.Range("BM10").Value = 0.867
$ws.Range("BN10").Value = 0.116
$ws.Range("BO10").Value = 0.34
$ws.Range("BP10").Value = 0.652
$ws.Range("BQ10").Value = 0.673
$ws.Range("E11").Value = 0.644
$ws.Range("F11").Value = 0.229
$ws.Range("G11").Value = 0.479
$ws.Range("N11").Value = 0.889
$ws.Range("O11").Value = 0.099
$ws.Range("P11").Value = 0.314
$ws.Range("W11").Value = 0.267
$ws.Range("X11").Value = 0.196
$ws.Range("Y11").Value = 0.442
$ws.Range("AI11").Value = 0.267
$ws.Range("AJ11").Value = 0.196
$ws.Range("AK11").Value = 0.442
$ws.Range("AU11").Value = 0.333
$ws.Range("AV11").Value = 0.222
$ws.Range("AW11").Value = 0.471
$ws.Range("BA11").Value = 1.956
$ws.Range("BB11").Value = 0.249
$ws.Range("BC11").Value = 0.499
$ws.Range("BG11").Value = 0.622
$ws.Range("BH11").Value = 0.235
$ws.Range("BI11").Value = 0.485
$ws.Range("BM11").Value = 0.867
$ws.Range("BN11").Value = 0.116
$ws.Range("BO11").Value = 0.34
$ws.Range("BP11").Value = 0.652
$ws.Range("BQ11").Value = 0.673
$ws.Range("E12").Value = 1.414
$ws.Range("F12").Value = 0.656
$ws.Range("G12").Value = 0.8100000000000001
$ws.Range("N12").Value = 1.488
$ws.Range("O12").Value = 1.03
$ws.Range("P12").Value = 1.015
$ws.Range("W12").Value = 1.75
$ws.Range("X12").Value = 0.6879999999999999
$ws.Range("Y12").Value = 0.829
$ws.Range("AI12").Value = 1.917
$ws.Range("AJ12").Value = 0.91
$ws.Range("AK12").Value = 0.954
$ws.Range("AU12").Value = 2.667
$ws.Range("AV12").Value = 1.689
$ws.Range("AW12").Value = 1.3
$ws.Range("BA12").Value = 3.694
$ws.Range("BB12").Value = 0.372
$ws.Range("BC12").Value = 0.61
$ws.Range("BG12").Value = 1.071
$ws.Range("BH12").Value = 0.066
$ws.Range("BI12").Value = 0.258
$ws.Range("BM12").Value = 1.385
$ws.Range("BN12").Value = 0.442
$ws.Range("BO12").Value = 0.665
$ws.Range("BP12").Value = 1.231
$ws.Range("BQ12").Value = 1.288
$ws.Range("E13").Value = 1.683
$ws.Range("F13").Value = 0.9419999999999999
$ws.Range("G13").Value = 0.971
$ws.Range("N13").Value = 2.287
$ws.Range("O13").Value = 1.252
$ws.Range("P13").Value = 1.119
$ws.Range("W13").Value = 1.088
$ws.Range("X13").Value = 0.176
$ws.Range("Y13").Value = 0.42
$ws.Range("AI13").Value = 1.383
$ws.Range("AJ13").Value = 0.401
$ws.Range("AK13").Value = 0.633
$ws.Range("AU13").Value = 2.487
$ws.Range("AV13").Value = 1.432
$ws.Range("AW13").Value = 1.197
$ws.Range("BA13").Value = 2.571
$ws.Range("BB13").Value = 0.332
$ws.Range("BC13").Value = 0.576
$ws.Range("BG13").Value = 0.639
$ws.Range("BH13").Value = 0.093
$ws.Range("BI13").Value = 0.305
$ws.Range("BM13").Value = 0.999
$ws.Range("BN13").Value = 0.392
$ws.Range("BO13").Value = 0.626
$ws.Range("BP13").Value = 0.857
$ws.Range("BQ13").Value = 0.798
